# Exporting metrics for each file in the analysed project.
#
# This recreates, via the Excel object model, the same change the commit
# made to the issues-template.xlsx workbook:
#   1. A new "Metrics" worksheet is added (after "Unconfirmed"), cloning the
#      layout of the existing "All" sheet (single "Colonne1" column table),
#      and becomes the active/selected sheet.
#   2. A matching "metrics" Excel table (ListObject) is added on that sheet.
#   3. The "selected" table on the "Issues" sheet already has a 10th column
#      ("Comments") that the pivot table/cache had not picked up yet; the
#      pivot cache + "synthesis" pivot table are rebuilt so the pivot field
#      list includes "Comments" too (cacheFields/pivotFields count 9 -> 10),
#      while keeping the exact same row-field layout and data field.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rebuild the pivot table/cache so the new "Comments" column of the
#    "selected" table is included as a pivot field.
# ---------------------------------------------------------------------
$wsIssues = $wb.Worksheets.Item("Issues")
$selectedTable = $wsIssues.ListObjects.Item(1)

$wsTCD = $wb.Worksheets.Item("TCD")
$oldPivot = $wsTCD.PivotTables(1)
$oldPivot.TableRange2.Delete()

$newCache = $wb.PivotCaches().Create(1, $selectedTable.Range)
$pivot = $newCache.CreatePivotTable($wsTCD.Range("A3"), "synthesis")

# Recreate the original row-field layout: Language, Type, Severity, Rule,
# Message (outer to inner), counting "Message".
$pivot.PivotFields("Language").Orientation = 1
$pivot.PivotFields("Type").Orientation = 1
$pivot.PivotFields("Severity").Orientation = 1
$pivot.PivotFields("Rule").Orientation = 1
$pivot.PivotFields("Message").Orientation = 1

$pivot.AddDataField($pivot.PivotFields("Message"), "Nombre de Message", -4112)

# CreatePivotTable defaults to using the last column ("Comments") as a
# Count data field - drop it again so only "Nombre de Message" remains,
# matching the original template.
$pivot.PivotFields("Comments").Orientation = 0

# ---------------------------------------------------------------------
# 2. Add the "Metrics" worksheet, cloning the "All" sheet's layout.
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("All")

$wsMetrics = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsMetrics.Name = "Metrics"

foreach ($col in @(2, 3, 10, 13, 15, 16)) {
    $wsMetrics.Columns.Item($col).ColumnWidth = $wsAll.Columns.Item($col).ColumnWidth
}

$headerCell = $wsMetrics.Range("A1")
$headerCell.Value = "Colonne1"
$headerCell.HorizontalAlignment = -4108
$headerCell.VerticalAlignment = -4108

$wsMetrics.PageSetup.PaperSize = 9
$wsMetrics.PageSetup.Orientation = 1

$metricsTable = $wsMetrics.ListObjects.Add(1, $wsMetrics.Range("A1:A2"), $null, 1)
$metricsTable.Name = "metrics"
$metricsTable.TableStyle = "TableStyleLight16"

# Make the new sheet the active / selected tab, like the original commit.
$wsMetrics.Activate()
$wsMetrics.Range("A1").Select()
